# Natmi Slit2-Robo4 LR-pair table, following Dr Hou's advice:
# expand the Sending/Target cluster combinations to the full FAPs/sCs/ECs
# cross-product (previously FAPs/sCs only on the sending side and the
# receptor/target side was missing ECs vs FAPs/sCs/ECs combinations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "sCs")

# Ligand-expressing-cell stats keyed by sending cluster
# (Ligand-expressing cells, detection rate, avg expr, total expr, avg-specificity, total-specificity)
$ligand = @{
  "ECs"  = @(2.0, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662)
  "FAPs" = @(3.0, 1.0, 0.37963, 1.13889, 0.1990797397484571, 0.1990797397484571)
  "sCs"  = @(3.0, 1.0, 1.493861, 4.481583000000001, 0.7833876645691064, 0.7833876645691064)
}

# Receptor-expressing-cell stats keyed by target cluster
# (Receptor-expressing cells, detection rate, avg expr, total expr, avg-specificity, total-specificity)
$receptor = @{
  "ECs"  = @(3.0, 1.0, 25.33077566666667, 75.992327, 0.988229338287255, 0.988229338287255)
  "FAPs" = @(1.0, 0.3333333333333333, 0.05766533333333334, 0.172996, 0.002249697164903793, 0.002249697164903793)
  "sCs"  = @(3.0, 1.0, 0.244046, 0.732138, 0.009520964547841182, 0.009520964547841182)
}

# Edge weight/specificity stats keyed by "sending|target"
# (avg weight, total weight, avg-specificity, total-specificity)
$edge = @{
  "ECs|ECs"   = @(0.8468922664555556, 7.622030398100001, 0.01732622542971233, 0.01732622542971233)
  "ECs|FAPs"  = @(0.001927944311111111, 0.0173514988, 0.00003944303080018216, 0.00003944303080018214)
  "ECs|sCs"   = @(0.008159271266666667, 0.0734334414, 0.0001669272219241125, 0.0001669272219241124)
  "FAPs|ECs"  = @(9.616322366336666, 86.54690129703, 0.1967364394780166, 0.1967364394780166)
  "FAPs|FAPs" = @(0.02189149049333333, 0.19702341444, 0.0004478691261018888, 0.0004478691261018888)
  "FAPs|sCs"  = @(0.09264718297999999, 0.8338246468199999, 0.001895431144338509, 0.001895431144338509)
  "sCs|ECs"   = @(37.84065786818234, 340.5659208136411, 0.774166673379526, 0.774166673379526)
  "sCs|FAPs"  = @(0.08614399251866668, 0.7752959326680001, 0.001762385008001723, 0.001762385008001722)
  "sCs|sCs"   = @(0.364570801606, 3.281137214454, 0.007458606181578562, 0.007458606181578562)
}

$row = 2
foreach ($sendingCluster in $clusters) {
  foreach ($targetCluster in $clusters) {
    $ws.Cells.Item($row, 1).Value = $sendingCluster
    $ws.Cells.Item($row, 2).Value = "Slit2"
    $ws.Cells.Item($row, 3).Value = "Robo4"
    $ws.Cells.Item($row, 4).Value = $targetCluster

    $l = $ligand[$sendingCluster]
    $ws.Cells.Item($row, 5).Value  = $l[0]
    $ws.Cells.Item($row, 6).Value  = $l[1]
    $ws.Cells.Item($row, 7).Value  = $l[2]
    $ws.Cells.Item($row, 8).Value  = $l[3]
    $ws.Cells.Item($row, 9).Value  = $l[4]
    $ws.Cells.Item($row, 10).Value = $l[5]

    $r = $receptor[$targetCluster]
    $ws.Cells.Item($row, 11).Value = $r[0]
    $ws.Cells.Item($row, 12).Value = $r[1]
    $ws.Cells.Item($row, 13).Value = $r[2]
    $ws.Cells.Item($row, 14).Value = $r[3]
    $ws.Cells.Item($row, 15).Value = $r[4]
    $ws.Cells.Item($row, 16).Value = $r[5]

    $e = $edge["$sendingCluster|$targetCluster"]
    $ws.Cells.Item($row, 17).Value = $e[0]
    $ws.Cells.Item($row, 18).Value = $e[1]
    $ws.Cells.Item($row, 19).Value = $e[2]
    $ws.Cells.Item($row, 20).Value = $e[3]

    $row = $row + 1
  }
}
